$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 28 - "estLeCreateur": move "X" from Terminé (D) to A faire (B), clear Responsable (E) and Date de fin (G)
$ws.Range("B28").Value = "X"
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = ""
# Clear the date and drop its date number-format (style goes back to the plain bordered style)
$ws.Range("F28").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = ""
$excel.CutCopyMode = $false

# Row 29 - "existeReservation": move "X" from A faire (B) to En cours (C), set Responsable (E) to "Ewan"
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = "X"
$ws.Range("E29").Value = "Ewan"

# Row 35 - "getReservation": clear Responsable (E)
$ws.Range("E35").Value = ""

# Row 36 - "getUtilisateur": move "X" from Terminé (D) to A faire (B)
$ws.Range("B36").Value = "X"
$ws.Range("D36").Value = ""

# Update the selection to reflect the saved cursor position
$ws.Range("F25").Select()
